$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: 11-Sep-2023
$ws.Range("A36").Value = [DateTime]"2023-09-11"
$ws.Range("B36").Value = "PRESENT"
$ws.Range("C36").Value = "PRESENT"
$ws.Range("D36").Value = "PRESENT"
$ws.Range("E36").Value = "PRESENT"
$ws.Range("F36").Value = "PRESENT"
$ws.Range("G36").Value = "ABSENT"
$ws.Range("H36").Value = "ABSENT"
$ws.Range("I36").Value = "PRESENT"
$ws.Range("J36").Value = "ABSENT"
$ws.Range("K36").Value = "ABSENT"

# Row 37: 12-Sep-2023
$ws.Range("A37").Value = [DateTime]"2023-09-12"
$ws.Range("B37").Value = "PRESENT"
$ws.Range("C37").Value = "PRESENT"
$ws.Range("D37").Value = "PRESENT"
$ws.Range("E37").Value = "PRESENT"
$ws.Range("F37").Value = "PRESENT"
$ws.Range("G37").Value = "PRESENT"
$ws.Range("H37").Value = "ABSENT"
$ws.Range("I37").Value = "PRESENT"
$ws.Range("J37").Value = "PRESENT"
$ws.Range("K37").Value = "ABSENT"

$ws.Range("A36:A37").NumberFormat = $ws.Range("A35").NumberFormat

$ws.Range("K37").Select()
